$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# The five group header rows (A1, A7, A13, A19, A25) used to be labelled with
# actual colour names ("Rød", "Gul", "Hvid", "Blå", "Grøn"). Rename them to the
# generic "Kategori N" labels used going forward.
$ws.Range("A1").Value = "Kategori 1"
$ws.Range("A7").Value = "Kategori 2"
$ws.Range("A13").Value = "Kategori 3"
$ws.Range("A19").Value = "Kategori 4"
$ws.Range("A25").Value = "Kategori 5"

# Reset the saved view of the sheet: no scrolled/frozen top-left cell, and the
# selection back at A2 instead of H30.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A2").Select()
